$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '56.060.46'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +6.90%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.506.68'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +7.29%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '491.80'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +12.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '141.90'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +15.75%  '
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +10.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.520.23'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +7.89%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0996'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +11.92%  '
$c.Style = 'Normal'

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.57'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +6.38%  '
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +8.40%  '
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +2.12%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.938.96'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +7.53%  '
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '55.966.47'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +6.73%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '20.99'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +10.21%  '
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +16.51%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.515.84'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +6.68%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.43'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +10.40%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '324.30'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +7.41%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.03'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +10.78%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.80'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +11.43%  '
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '58.45'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +7.92%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.171'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +12.42%  '
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +11.99%  '
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -6.27%  '
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.616.02'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +7.48%  '
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.51'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +5.88%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0₃0805'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +18.23%  '
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.22%  '
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '150.72'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +4.41%  '
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '18.37'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +7.51%  '
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +12.52%  '
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +8.74%  '
$c.Style = 'Normal'

$c = $ws.Range('B36')
$c.NumberFormat = '@'
$c.Value = 'NEARProtocol'
$c.Style = 'Normal'
$c = $ws.Range('C36')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.75'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +6.93%  '
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +12.76%  '
$c.Style = 'Normal'

$c = $ws.Range('B38')
$c.NumberFormat = '@'
$c.Value = 'Fetch.AI'
$c.Style = 'Normal'
$c = $ws.Range('C38')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.880'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +5.49%  '
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '34.50'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +7.84%  '
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.617'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +16.50%  '
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +9.12%  '
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +0.25%  '
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +9.15%  '
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.34'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +9.35%  '
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '4.80'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +16.17%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.017.65'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +4.47%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '260.45'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +34.82%  '
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +10.57%  '
$c.Style = 'Normal'

$c = $ws.Range('B49')
$c.NumberFormat = '@'
$c.Value = 'VeChain'
$c.Style = 'Normal'
$c = $ws.Range('C49')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0227'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +8.34%  '
$c.Style = 'Normal'

$c = $ws.Range('B50')
$c.NumberFormat = '@'
$c.Value = 'WhiteBITCoin'
$c.Style = 'Normal'
$c = $ws.Range('C50')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '10.14'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '17.65'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +11.76%  '
$c.Style = 'Normal'
